# Add a new LeetCode tracking row (#2040 "Kth Smallest Product of Two Sorted
# Arrays") to the bottom of the tracking table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (13) down into the new row
# (14) first. Using PasteSpecial(xlPasteFormats) re-uses the existing style
# records (general/center, wrap-center, date) instead of minting brand new
# ones, which keeps styles.xml untouched - exactly like the source row.
$ws.Range("A13:I13").Copy()
$ws.Range("A14:I14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 34

# Fill in the new problem's data.
$ws.Range("A14").Value = 2040
$ws.Range("B14").Value = "Kth Smallest Product of Two Sorted Arrays"
$ws.Range("C14").Value = "#array #binary-search "
$ws.Range("D14").Value = "hard"
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 60
$ws.Range("H14").Value = 45833
$ws.Range("I14").Value = 45833

# Move the sheet's active selection the way the author left it.
[void]$ws.Range("L13").Select()
